$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16, column A: rename "Customer" task to "Member"
$ws.Range("A16").Value = "Member"

# New column E: mark completion status (1) for the relevant rows
$ws.Range("E2").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("E15").Value = 1
$ws.Range("E16").Value = 1

# Update selection to match the author's last selected cell
$ws.Range("E6").Select()
